# Introduce the "addTag" privilege rule for Task objects.
# Inserts a new rule row ("Task – Only participants can add tags") right
# after the existing "Task – anybody can subscribe" row (row 53) and before
# the "Folder – default public access" row, pushing the Folder rows down by
# one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new blank row at row 54 (old rows 54-55 shift down to 55-56).
$ws.Rows("54:54").Insert() | Out-Null

# Copy formatting from the row above (row 53, a same-shape TASK rule row)
# so the new row reuses the existing cell styles instead of creating new
# style entries.
$ws.Range("A53:G53").Copy() | Out-Null
$ws.Range("A54:G54").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Fill in the new rule's content.
$ws.Cells.Item(54, 2).Value2 = "Task – Only participants can add tags"
$ws.Cells.Item(54, 3).Value2 = "TASK"
$ws.Cells.Item(54, 7).Value2 = "grant addTag to assignee, co-owner, supervisor, owning group, approver, collaborator, reader"

# The new rule text wraps onto 3 lines, so the row is taller (45pt) than the
# single-line "anybody can subscribe" row it was copied from (30pt).
$ws.Rows("54:54").RowHeight = 45

# Reflect the cursor/scroll position left behind by the edit.
[void]$excel.Goto($ws.Range("A46"), $true)
$ws.Range("C54").Select() | Out-Null
